# compliance_matrix.xlsx update
# - Start Date (col G) 2025-01-01 -> 2025-07-01 (45658 -> 45839) for all data rows
# - End Date   (col H) 2025-09-25 -> 2025-09-30 (45925 -> 45930) for all data rows
# - Email (col D) for rows 15-21 and 26: yashalifarooqui30@gmail.com -> farooquiyashal@gmail.com
# - Rebuild hyperlinks so the D/E mailto + sharepoint links stay correct
# - Update the active selection to G3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmail = "farooquiyashal@gmail.com"
$oldEmail = "yashalifarooqui30@gmail.com"
$sharepointLink = "https://nflpk.sharepoint.com/nfl_site/IT%20Department/Forms/AllItems.aspx?id=%2Fnfl%5Fsite%2FIT%20Department%2FIT%2DGovernance%2FCOMPLIANCE%20MATRIX%2FEVIDENCE%2FApplication%20Nfoods%2Ecom&viewid=15a9448d%2D5114%2D4f3f%2D81c3%2Ddaf0476b38d2&CT=1756293467065&OR=OWA%2DNT%2DMail&CID=bedfc76e%2D9e1e%2Db098%2D9316%2D9bdf1d6e274a&csf=1&web=1&e=svmXMI&FolderCTID=0x01200095CBBE0C68C0474BAF872F72573A9D97"

# 1. Update Start Date / End Date for every data row (2-29)
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 7).Value = 45839
    $ws.Cells.Item($r, 8).Value = 45930
}

# 2. Update the email text for the rows whose contact changed
$emailRows = @(15, 16, 17, 18, 19, 20, 21, 26)
foreach ($r in $emailRows) {
    $ws.Cells.Item($r, 4).Value = $newEmail
}

# 3. Rebuild the hyperlinks collection in the order Excel would leave it:
#    existing links that did not change keep their relative order, the
#    (deleted + re-added) D-column links land at the end.
$ws.Hyperlinks.Delete()

$unchangedMailtoRows = @(
    @(3,  "farooquiyashal@gmail.com"),
    @(2,  "farooquiyashal@gmail.com"),
    @(4,  "farooquiyashal@gmail.com"),
    @(5,  "farooquiyashal@gmail.com"),
    @(6,  "farooquiyashal@gmail.com"),
    @(7,  "aliyashal309@gmail.com"),
    @(8,  "aliyashal309@gmail.com"),
    @(9,  "aliyashal309@gmail.com"),
    @(10, "aliyashal309@gmail.com"),
    @(11, "aliyashal309@gmail.com"),
    @(12, "aliyashal309@gmail.com"),
    @(13, "aliyashal309@gmail.com"),
    @(14, "aliyashal309@gmail.com"),
    @(22, "yashal.ali@nfoods.com"),
    @(23, "yashal.ali@nfoods.com"),
    @(24, "yashal.ali@nfoods.com"),
    @(25, "yashal.ali@nfoods.com"),
    @(27, "aliyashal309@gmail.com"),
    @(28, "aliyashal309@gmail.com"),
    @(29, "aliyashal309@gmail.com")
)

foreach ($pair in $unchangedMailtoRows) {
    $row = $pair[0]
    $addr = $pair[1]
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 4), "mailto:$addr")
}

# Attachment-link (column E) hyperlinks
$ws.Hyperlinks.Add($ws.Range("E3"), $sharepointLink)
$ws.Hyperlinks.Add($ws.Range("E4:E29"), $sharepointLink, "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("E2"), $sharepointLink)

# The D-column links whose target email changed, re-created last
$changedMailtoRows = @(15, 16, 17, 18, 19, 20, 21, 26)
foreach ($r in $changedMailtoRows) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 4), "mailto:$newEmail")
}

# 4. Move the selection to G3 to match the saved view state
$ws.Range("G3").Select()

Write-Host "edit complete"
